# Applies the "Finished burgers optimality. Needs polish" edit:
#  - Adds two new gradebook columns (I: "Q3", J: "W3") with per-student scores
#  - Updates the window view / selection metadata

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: match the saved window position/size from the edited workbook.
# (Some hosts may not persist on-screen window chrome metrics into the saved
# file, but setting them here is harmless if unsupported.)
try {
    $win = $excel.ActiveWindow
    $win.Left = 720
    $win.Top = 720
    $win.Width = 19740
    $win.Height = 14000
} catch {
}

# New header cells for the two added columns.
$ws.Range("I1").Value = "Q3"
$ws.Range("J1").Value = "W3"

# Per-row data for the new columns (row => I/J values), taken from the diff.
$data = @(
    @{Row=2; I=4; J=49},
    @{Row=3; I=0; J=0},
    @{Row=4; I=4; J=50},
    @{Row=5; I=4; J=50},
    @{Row=6; I=4; J=50},
    @{Row=7; I=4; J=50},
    @{Row=8; I=4; J=50},
    @{Row=9; I=4; J=50},
    @{Row=10; I=4; J=50},
    @{Row=11; I=4; J=50},
    @{Row=12; I=4; J=50},
    @{Row=13; I=4; J=50},
    @{Row=14; I=2; J=50},
    @{Row=15; I=4; J=50},
    @{Row=16; I=0; J=0},
    @{Row=17; I=4; J=50},
    @{Row=18; I=4; J=50},
    @{Row=19; I=4; J=50},
    @{Row=20; I=4; J=50},
    @{Row=21; I=2; J=50},
    @{Row=22; I=4; J=49},
    @{Row=23; I=4; J=50},
    @{Row=24; I=4; J=50},
    @{Row=25; I=4; J=50},
    @{Row=26; I=4; J=50},
    @{Row=27; I=4; J=50},
    @{Row=28; I=4; J=50},
    @{Row=29; I=4; J=49},
    @{Row=30; I=4; J=50},
    @{Row=31; I=4; J=50},
    @{Row=32; I=4; J=50},
    @{Row=33; I=4; J=50}
)

foreach ($d in $data) {
    $ws.Cells.Item($d.Row, 9).Value = $d.I
    $ws.Cells.Item($d.Row, 10).Value = $d.J
}

# Match the saved selection/active cell from the edited workbook.
$ws.Range("J16").Select()
